# Insert a new "Source" column between "Species" (A) and "Description" (B).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column B; this shifts the existing
# "Description" column (and its formatting) from B to C.
$ws.Columns("B:B").Insert()

# Determine the last used row (should be 45: header + 44 species rows).
$lastRow = $ws.Cells(1, 1).End(-4121).Row  # xlDown = -4121

# Header cell, mirrors the style/formatting used by the other header cells.
$ws.Cells.Item(1, 2).Value = "Source"
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)  # xlPasteFormats

# Fill the "Source" value for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "Inventario IEET - Peces"
}
